$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Lombardia"
$ws.Range("C2").Value = 4140912
$ws.Range("B3").Value = "Veneto"
$ws.Range("C3").Value = 2710818
$ws.Range("B4").Value = "Campania"
$ws.Range("C4").Value = 2454926
$ws.Range("B5").Value = "Lazio"
$ws.Range("C5").Value = 2401611
$ws.Range("B6").Value = "Emilia-Romagna"
$ws.Range("C6").Value = 2145422
$ws.Range("B7").Value = "Sicilia"
$ws.Range("C7").Value = 1822575
$ws.Range("B8").Value = "Piemonte"
$ws.Range("C8").Value = 1727244
$ws.Range("B9").Value = "Puglia"
$ws.Range("C9").Value = 1632111
$ws.Range("B10").Value = "Toscana"
$ws.Range("C10").Value = 1596603
$ws.Range("B11").Value = "Marche"
$ws.Range("C11").Value = 715924
$ws.Range("B12").Value = "Liguria"
$ws.Range("C12").Value = 664018
$ws.Range("B13").Value = "Abruzzo"
$ws.Range("C13").Value = 654796
$ws.Range("B14").Value = "Calabria"
$ws.Range("C14").Value = 634195
$ws.Range("B15").Value = "Friuli Venezia Giulia"
$ws.Range("C15").Value = 577967
$ws.Range("B16").Value = "Sardegna"
$ws.Range("C16").Value = 511980
$ws.Range("B17").Value = "Umbria"
$ws.Range("C17").Value = 440864
$ws.Range("B18").Value = "P.A. Bolzano"
$ws.Range("C18").Value = 295395
$ws.Range("B19").Value = "P.A. Trento"
$ws.Range("C19").Value = 245353
$ws.Range("B20").Value = "Basilicata"
$ws.Range("C20").Value = 200103
$ws.Range("B21").Value = "Molise"
$ws.Range("C21").Value = 102127
$ws.Range("B22").Value = "Valle d'Aosta"
$ws.Range("C22").Value = 50647
